# "add shunt and generator" - adds two new columns (bsh, options) to the
# "node" worksheet, with a shunt susceptance value in row 2 and a
# generator "base" marker in row 2 / column "options".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("node")

# New headers.
# Shared-string insertion order matters for reproducing the same string
# table as the target workbook: "options" must be registered before
# "bsh" so that they end up as shared-string indices 18 and 19
# respectively, and "base" is registered afterwards as index 20.
$ws.Range("N1").Value = "options"
$ws.Range("M1").Value = "bsh"

# Row 2 (bus 1): shunt susceptance 0, and an "options" tag of "base"
# (this is the new generator / shunt reference row).
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = "base"

# Remaining rows only get the new "bsh" shunt susceptance column,
# defaulting to 0 (no shunt attached).
$ws.Range("M3").Value = 0
$ws.Range("M4").Value = 0
$ws.Range("M5").Value = 0
$ws.Range("M6").Value = 0
$ws.Range("M7").Value = 0
$ws.Range("M8").Value = 0

# Reproduce the final cell selection left behind on the "node" sheet.
$ws.Activate() | Out-Null
$ws.Range("O3").Select() | Out-Null
